$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Hide/zero-width the columns that now fall "inside" the existing
#     hidden helper-column block (J:AD instead of the previous J:W). ---
$ws.Range("X1:AD1").ColumnWidth = -0.9
$ws.Range("X1:AD1").EntireColumn.Hidden = $true

# --- Add a new column AJ ("27-jul") with the same look & feel as
#     the previous last column (AI). ---
$ws.Range("AI1:AI11").Copy()
$ws.Range("AJ1:AJ11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("AJ1").Value = "27-jul"
$ws.Range("AJ2").Value = 12
$ws.Range("AJ3").Value = 16
$ws.Range("AJ4").Value = 10
$ws.Range("AJ5").Value = 13
$ws.Range("AJ6").Value = 16
$ws.Range("AJ7").Value = 14
$ws.Range("AJ8").Value = 14
$ws.Range("AJ9").Value = 14
$ws.Range("AJ10").Value = 20
$ws.Range("AJ11").Value = 17

# --- Match the saved selection in the source file. ---
$null = $ws.Range("AQ17").Select()
